# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The forecast-error matrix stores, per row (per forecast horizon), a
# diagonal/triangular series of quarter-over-quarter naive errors. A newly
# computed (most-recent) error is inserted at the front of each row
# (column B); all previously stored errors shift one column to the right
# (B->C, C->D, ... J->K). Rows that were already fully populated out to
# column K drop their oldest value (previously in column K) when the new
# value is inserted. Rows that were not yet fully populated simply grow by
# one column and nothing is dropped. Row 16, which previously held no
# numeric values, receives its first value in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = New-Object 'object[,]' 1,10
$row[0,0] = -0.7322633397437844
$row[0,1] = 0.2879090979994584
$row[0,2] = -1.425880358899853
$row[0,3] = 1.402475014699119
$row[0,4] = -0.8328575851670005
$row[0,5] = 0.5184774727506619
$row[0,6] = -0.1029604570662399
$row[0,7] = 0.3976002401245912
$row[0,8] = -0.2703078322215502
$row[0,9] = 0.1586931430164528
$ws.Range("B2:K2").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0.2703549766394939
$row[0,1] = -1.443434480259818
$row[0,2] = 1.384920893339154
$row[0,3] = -0.8504117065269649
$row[0,4] = 0.5009233513906975
$row[0,5] = -0.1205145784262043
$row[0,6] = 0.3800461187646267
$row[0,7] = -0.2878619535815147
$row[0,8] = 0.1411390216564884
$row[0,9] = 0.1988105702346985
$ws.Range("B3:K3").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -1.355327161308811
$row[0,1] = 1.473028212290161
$row[0,2] = -0.7623043875759586
$row[0,3] = 0.5890306703417038
$row[0,4] = -0.0324072594751981
$row[0,5] = 0.4681534377156329
$row[0,6] = -0.1997546346305085
$row[0,7] = 0.2292463406074946
$row[0,8] = 0.2869178891857047
$row[0,9] = 0.4108842600239903
$ws.Range("B4:K4").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 1.651602845777944
$row[0,1] = -0.5837297540881751
$row[0,2] = 0.7676053038294873
$row[0,3] = 0.1461673740125855
$row[0,4] = 0.6467280712034165
$row[0,5] = -0.02118000114272489
$row[0,6] = 0.4078209740952782
$row[0,7] = 0.4654925226734883
$row[0,8] = 0.5894588935117738
$row[0,9] = -0.303959229723018
$ws.Range("B5:K5").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0.3282974736644749
$row[0,1] = 1.679632531582137
$row[0,2] = 1.058194601765235
$row[0,3] = 1.558755298956066
$row[0,4] = 0.8908472266099251
$row[0,5] = 1.319848201847928
$row[0,6] = 1.377519750426138
$row[0,7] = 1.501486121264424
$row[0,8] = 0.608067998029632
$row[0,9] = 1.378198724973394
$ws.Range("B6:K6").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0.7356582956163805
$row[0,1] = 0.1142203657994787
$row[0,2] = 0.6147810629903097
$row[0,3] = -0.0531270093558317
$row[0,4] = 0.3758739658821714
$row[0,5] = 0.4335455144603815
$row[0,6] = 0.557511885298667
$row[0,7] = -0.3359062379361248
$row[0,8] = 0.4342244890076376
$row[0,9] = 0.1683237681281231
$ws.Range("B7:K7").Value2 = $row

$row = New-Object 'object[,]' 1,9
$row[0,0] = 0.1181882633125878
$row[0,1] = 0.6187489605034189
$row[0,2] = -0.04915911184272259
$row[0,3] = 0.3798418633952805
$row[0,4] = 0.4375134119734906
$row[0,5] = 0.5614797828117761
$row[0,6] = -0.3319383404230157
$row[0,7] = 0.4381923865207467
$row[0,8] = 0.1722916656412322
$ws.Range("B8:J8").Value2 = $row

$row = New-Object 'object[,]' 1,8
$row[0,0] = 0.7543890506736601
$row[0,1] = 0.08648097832751878
$row[0,2] = 0.5154819535655218
$row[0,3] = 0.573153502143732
$row[0,4] = 0.6971198729820175
$row[0,5] = -0.1962982502527744
$row[0,6] = 0.5738324766909881
$row[0,7] = 0.3079317558114735
$ws.Range("B9:I9").Value2 = $row

$row = New-Object 'object[,]' 1,7
$row[0,0] = -0.1543252035281459
$row[0,1] = 0.2746757717098572
$row[0,2] = 0.3323473202880673
$row[0,3] = 0.4563136911263528
$row[0,4] = -0.4371044321084391
$row[0,5] = 0.3330262948353234
$row[0,6] = 0.06712557395580883
$ws.Range("B10:H10").Value2 = $row

$row = New-Object 'object[,]' 1,6
$row[0,0] = 0.2293445564577608
$row[0,1] = 0.2870161050359709
$row[0,2] = 0.4109824758742565
$row[0,3] = -0.4824356473605354
$row[0,4] = 0.287695079583227
$row[0,5] = 0.02179435870371246
$ws.Range("B11:G11").Value2 = $row

$row = New-Object 'object[,]' 1,5
$row[0,0] = 0.2201546830999171
$row[0,1] = 0.3441210539382026
$row[0,2] = -0.5492970692965893
$row[0,3] = 0.2208336576471732
$row[0,4] = -0.04506706323234141
$ws.Range("B12:F12").Value2 = $row

$row = New-Object 'object[,]' 1,4
$row[0,0] = 0.314534851581486
$row[0,1] = -0.5788832716533059
$row[0,2] = 0.1912474552904566
$row[0,3] = -0.07465326558905801
$ws.Range("B13:E13").Value2 = $row

$row = New-Object 'object[,]' 1,3
$row[0,0] = -0.5970339283829468
$row[0,1] = 0.1730967985608157
$row[0,2] = -0.0928039223186989
$ws.Range("B14:D14").Value2 = $row

$row = New-Object 'object[,]' 1,2
$row[0,0] = 0.1550649743121164
$row[0,1] = -0.1108357465673982
$ws.Range("B15:C15").Value2 = $row

$row = New-Object 'object[,]' 1,1
$row[0,0] = -0.1624199859130616
$ws.Range("B16:B16").Value2 = $row

